# Apply the "optimization study" parameter updates to MITHEx_inputs.xlsx
$wb = $excel.ActiveWorkbook

$wsPlant  = $wb.Worksheets.Item("Plant Description")
$wsHX     = $wb.Worksheets.Item("HX Parameters")
$wsCycle  = $wb.Worksheets.Item("Cycle Parameters")
$wsInput  = $wb.Worksheets.Item("Input options")

# --- Plant Description: remove the two now-unused "Intermediate" rows ---
$wsPlant.Rows.Item(13).Delete()
$wsPlant.Rows.Item(12).Delete()

# --- HX Parameters: update values to match the recent optimization studies ---
$wsHX.Range("B2").Value = 0.00135
$wsHX.Range("B4").Value = 100
$wsHX.Range("B5").Value = 150

# Switch the page to portrait orientation
$wsHX.PageSetup.Orientation = 1

# --- Restore/update selections on each sheet, finishing with HX Parameters active ---
$wsPlant.Activate()
$wsPlant.Range("B17").Select() | Out-Null

$wsInput.Activate()
$wsInput.Range("C36").Select() | Out-Null

$wsCycle.Activate()

$wsHX.Activate()
$wsHX.Range("C13").Select() | Out-Null

$wb.Save()
